$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix / add proper Indonesian descriptions for the "delete" API rows that
# previously shared the generic placeholder text "Menghapus Data "
$ws.Range("C42").Value = "Menghapus Data Akun E-Mail Orang"
$ws.Range("C43").Value = "Menghapus Data Akun Sosial Media Orang"
$ws.Range("C44").Value = "Menghapus Data Jenis Kelamin Orang"
$ws.Range("C45").Value = "Menghapus Data Produk"
$ws.Range("C46").Value = "Menghapus Data Jenis Produk"

# Fix typo ("EMail" -> "E-Mail") in the initialize Person Account EMail description
$ws.Range("C69").Value = "Menginisialisasi Data Akun E-Mail Orang"

# Update the stored selection in the sheet view to match the new active cell
$ws.Range("C47").Select()
